# Changes made for FinancePayout and InterpreAptTabs Test cases
#
# 1. "New appointment" sheet: K2/L2/K4/L4/K7/L7 "Abramson Building"/"CSE"
#    -> "Main Building"/"CS"; K6/L6 get the same values plus the border
#    style already used by the other K/L cells in that column; selection
#    moves to L6:L7 and this sheet becomes the active tab.
# 2. "Finance-Admin" sheet: selection moves to A22 and it is no longer the
#    active tab.
# 3. "Finance_Payout" sheet: B2/B3/B4 switch from
#    sravani.bandaru@sstech.us to ravi.thota@sstech.us, the B3:B4
#    hyperlink is removed, and the selection moves to C7.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Finance_Payout
# ---------------------------------------------------------------
$wsPayout = $wb.Worksheets.Item("Finance_Payout")

$wsPayout.Range("B2").Value = "ravi.thota@sstech.us"
$wsPayout.Range("B3").Value = "ravi.thota@sstech.us"
$wsPayout.Range("B4").Value = "ravi.thota@sstech.us"

# Remove the shared hyperlink covering B3:B4 (keep B2's own hyperlink).
for ($i = $wsPayout.Hyperlinks.Count(); $i -ge 1; $i--) {
    $link = $wsPayout.Hyperlinks.Item($i)
    $addr = $link.Range().Address()
    if ($addr -eq '$B$3:$B$4') {
        $link.Delete()
    }
}

$wsPayout.Range("C7").Select()

# ---------------------------------------------------------------
# Finance-Admin
# ---------------------------------------------------------------
$wsFinanceAdmin = $wb.Worksheets.Item("Finance-Admin")
$wsFinanceAdmin.Range("A22").Select()

# ---------------------------------------------------------------
# New appointment
# ---------------------------------------------------------------
$wsNewAppt = $wb.Worksheets.Item("New appointment")

$wsNewAppt.Range("K2").Value = "Main Building"
$wsNewAppt.Range("L2").Value = "CS"

$wsNewAppt.Range("K4").Value = "Main Building"
$wsNewAppt.Range("L4").Value = "CS"

$wsNewAppt.Range("K7").Value = "Main Building"
$wsNewAppt.Range("L7").Value = "CS"

# K6/L6 need the same border style already applied to K2/L2, so copy the
# formatting across before writing the new values.
$wsNewAppt.Range("K2").Copy()
$wsNewAppt.Range("K6").PasteSpecial(-4122)
$wsNewAppt.Range("L2").Copy()
$wsNewAppt.Range("L6").PasteSpecial(-4122)

$wsNewAppt.Range("K6").Value = "Main Building"
$wsNewAppt.Range("L6").Value = "CS"

# Make "New appointment" the active sheet/tab with L6:L7 selected - do
# this last so its tabSelected flag sticks.
$wsNewAppt.Activate()
$wsNewAppt.Range("L6:L7").Select()
